# Update NATMI LR-pair TPM-derived metrics (Fndc5-Itgb5) for rows 2-10
# to reflect the recomputed values with the new TPM inputs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.088045
$ws.Cells.Item(2, 8).Value = 0.264135
$ws.Cells.Item(2, 9).Value = 0.1003389266487061
$ws.Cells.Item(2, 10).Value = 0.1003389266487061
$ws.Cells.Item(2, 13).Value = 8.033114333333334
$ws.Cells.Item(2, 14).Value = 24.099343
$ws.Cells.Item(2, 15).Value = 0.1374088679258946
$ws.Cells.Item(2, 16).Value = 0.1374088679258946
$ws.Cells.Item(2, 17).Value = 0.7072755514783333
$ws.Cells.Item(2, 18).Value = 6.365479963305001
$ws.Cells.Item(2, 19).Value = 0.01378745831969807
$ws.Cells.Item(2, 20).Value = 0.01378745831969808

$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.088045
$ws.Cells.Item(3, 8).Value = 0.264135
$ws.Cells.Item(3, 9).Value = 0.1003389266487061
$ws.Cells.Item(3, 10).Value = 0.1003389266487061
$ws.Cells.Item(3, 15).Value = 0.6355200716780686
$ws.Cells.Item(3, 16).Value = 0.6355200716780686
$ws.Cells.Item(3, 17).Value = 3.271170310595
$ws.Cells.Item(3, 18).Value = 29.440532795355
$ws.Cells.Item(3, 19).Value = 0.06376740185588614
$ws.Cells.Item(3, 20).Value = 0.06376740185588614

$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.088045
$ws.Cells.Item(4, 8).Value = 0.264135
$ws.Cells.Item(4, 9).Value = 0.1003389266487061
$ws.Cells.Item(4, 10).Value = 0.1003389266487061
$ws.Cells.Item(4, 13).Value = 13.27489133333333
$ws.Cells.Item(4, 14).Value = 39.824674
$ws.Cells.Item(4, 15).Value = 0.2270710603960369
$ws.Cells.Item(4, 16).Value = 0.2270710603960369
$ws.Cells.Item(4, 17).Value = 1.168787807443333
$ws.Cells.Item(4, 18).Value = 10.51909026699
$ws.Cells.Item(4, 19).Value = 0.02278406647312185
$ws.Cells.Item(4, 20).Value = 0.02278406647312185

$ws.Cells.Item(5, 7).Value = 0.3888126666666667
$ws.Cells.Item(5, 9).Value = 0.4431034770941504
$ws.Cells.Item(5, 10).Value = 0.4431034770941504
$ws.Cells.Item(5, 13).Value = 8.033114333333334
$ws.Cells.Item(5, 14).Value = 24.099343
$ws.Cells.Item(5, 15).Value = 0.1374088679258946
$ws.Cells.Item(5, 16).Value = 0.1374088679258946
$ws.Cells.Item(5, 17).Value = 3.123376605581556
$ws.Cells.Item(5, 18).Value = 28.110389450234
$ws.Cells.Item(5, 19).Value = 0.06088634716153476
$ws.Cells.Item(5, 20).Value = 0.06088634716153477

$ws.Cells.Item(6, 7).Value = 0.3888126666666667
$ws.Cells.Item(6, 9).Value = 0.4431034770941504
$ws.Cells.Item(6, 10).Value = 0.4431034770941504
$ws.Cells.Item(6, 15).Value = 0.6355200716780686
$ws.Cells.Item(6, 16).Value = 0.6355200716780686
$ws.Cells.Item(6, 19).Value = 0.2816011535236759
$ws.Cells.Item(6, 20).Value = 0.2816011535236759

$ws.Cells.Item(7, 7).Value = 0.3888126666666667
$ws.Cells.Item(7, 9).Value = 0.4431034770941504
$ws.Cells.Item(7, 10).Value = 0.4431034770941504
$ws.Cells.Item(7, 13).Value = 13.27489133333333
$ws.Cells.Item(7, 14).Value = 39.824674
$ws.Cells.Item(7, 15).Value = 0.2270710603960369
$ws.Cells.Item(7, 16).Value = 0.2270710603960369
$ws.Cells.Item(7, 17).Value = 5.161445899023557
$ws.Cells.Item(7, 18).Value = 46.45301309121201
$ws.Cells.Item(7, 19).Value = 0.1006159764089398
$ws.Cells.Item(7, 20).Value = 0.1006159764089398

$ws.Cells.Item(8, 7).Value = 0.4006183333333334
$ws.Cells.Item(8, 8).Value = 1.201855
$ws.Cells.Item(8, 9).Value = 0.4565575962571436
$ws.Cells.Item(8, 10).Value = 0.4565575962571436
$ws.Cells.Item(8, 13).Value = 8.033114333333334
$ws.Cells.Item(8, 14).Value = 24.099343
$ws.Cells.Item(8, 15).Value = 0.1374088679258946
$ws.Cells.Item(8, 16).Value = 0.1374088679258946
$ws.Cells.Item(8, 17).Value = 3.218212875696111
$ws.Cells.Item(8, 18).Value = 28.963915881265
$ws.Cells.Item(8, 19).Value = 0.06273506244466175
$ws.Cells.Item(8, 20).Value = 0.06273506244466176

$ws.Cells.Item(9, 7).Value = 0.4006183333333334
$ws.Cells.Item(9, 8).Value = 1.201855
$ws.Cells.Item(9, 9).Value = 0.4565575962571436
$ws.Cells.Item(9, 10).Value = 0.4565575962571436
$ws.Cells.Item(9, 15).Value = 0.6355200716780686
$ws.Cells.Item(9, 16).Value = 0.6355200716780686
$ws.Cells.Item(9, 17).Value = 14.88432958010167
$ws.Cells.Item(9, 18).Value = 133.958966220915
$ws.Cells.Item(9, 19).Value = 0.2901515162985066
$ws.Cells.Item(9, 20).Value = 0.2901515162985066

$ws.Cells.Item(10, 7).Value = 0.4006183333333334
$ws.Cells.Item(10, 8).Value = 1.201855
$ws.Cells.Item(10, 9).Value = 0.4565575962571436
$ws.Cells.Item(10, 10).Value = 0.4565575962571436
$ws.Cells.Item(10, 13).Value = 13.27489133333333
$ws.Cells.Item(10, 14).Value = 39.824674
$ws.Cells.Item(10, 15).Value = 0.2270710603960369
$ws.Cells.Item(10, 16).Value = 0.2270710603960369
$ws.Cells.Item(10, 17).Value = 5.318164841141112
$ws.Cells.Item(10, 18).Value = 47.86348357027001
$ws.Cells.Item(10, 19).Value = 0.1036710175139753
$ws.Cells.Item(10, 20).Value = 0.1036710175139753
